$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 232.61539
$ws.Range("I6").Value = 152.4
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 457.2
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -345.2
$ws.Range("N6").Value = -1724
# Row 8
$ws.Range("H8").Value = 16.2
$ws.Range("I8").Value = 16.2
$ws.Range("K8").Value = 48.59999999999999
$ws.Range("M8").Value = 90.40000000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 3019.1667
$ws.Range("I74").Value = 3019.1667
$ws.Range("K74").Value = 3019.1667
$ws.Range("M74").Value = -2145.1667
# Row 77
$ws.Range("H77").Value = 3019.1667
$ws.Range("I77").Value = 3019.1667
$ws.Range("K77").Value = 15095.8335
$ws.Range("M77").Value = -10727.8335
# Row 122
$ws.Range("H122").Value = 2621.5
$ws.Range("I122").Value = 2026.875
$ws.Range("K122").Value = 6080.625
$ws.Range("M122").Value = -3630.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 4427.6665
$ws.Range("I36").Value = 2934.2
$ws.Range("K36").Value = 2934.2
$ws.Range("M36").Value = -2400.2
# Row 99
$ws.Range("H99").Value = 4957.8
$ws.Range("I99").Value = 5447.25
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 5447.25
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -3949.25
$ws.Range("N99").Value = -5996
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4644.2
$ws.Range("I31").Value = 3930.25
$ws.Range("K31").Value = 3930.25
$ws.Range("M31").Value = -3635.25
# Row 34
$ws.Range("H34").Value = 4644.2
$ws.Range("I34").Value = 3930.25
$ws.Range("K34").Value = 3930.25
$ws.Range("M34").Value = -3728.25
# Row 58
$ws.Range("H58").Value = 2704.5625
$ws.Range("I58").Value = 2527.4167
$ws.Range("J58").Value = 3236
$ws.Range("K58").Value = 2527.4167
$ws.Range("L58").Value = 3236
$ws.Range("M58").Value = -2324.4167
$ws.Range("N58").Value = -3642
# Row 64
$ws.Range("H64").Value = 40000
$ws.Range("I64").Value = 30000
$ws.Range("K64").Value = 30000
$ws.Range("M64").Value = -29752
# Row 67
$ws.Range("H67").Value = 40000
$ws.Range("I67").Value = 30000
$ws.Range("K67").Value = 30000
$ws.Range("M67").Value = -29142
# Row 99
$ws.Range("H99").Value = 8659.799999999999
$ws.Range("I99").Value = 8659.799999999999
$ws.Range("K99").Value = 8659.799999999999
$ws.Range("M99").Value = -7161.799999999999
# Row 105
$ws.Range("H105").Value = 1750
$ws.Range("I105").Value = 1625
$ws.Range("K105").Value = 1625
$ws.Range("M105").Value = 122
# Row 126
$ws.Range("H126").Value = 8659.799999999999
$ws.Range("I126").Value = 8659.799999999999
$ws.Range("K126").Value = 25979.4
$ws.Range("M126").Value = -23509.4
# Row 136
$ws.Range("H136").Value = 2704.5625
$ws.Range("I136").Value = 2527.4167
$ws.Range("J136").Value = 3236
$ws.Range("K136").Value = 7582.250100000001
$ws.Range("L136").Value = 9708
$ws.Range("M136").Value = -5032.250100000001
$ws.Range("N136").Value = -14808

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 7720.5713
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 6000
$ws.Range("M81").Value = -4877
# Row 84
$ws.Range("H84").Value = 7720.5713
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 18000
$ws.Range("M84").Value = -12384
# Row 121
$ws.Range("H121").Value = 875.6
$ws.Range("I121").Value = 496
$ws.Range("J121").Value = 970.5
$ws.Range("K121").Value = 1488
$ws.Range("L121").Value = 2911.5
$ws.Range("M121").Value = -178
$ws.Range("N121").Value = -5531.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -724
# Row 80
$ws.Range("H80").Value = 3610.9285
$ws.Range("I80").Value = 2616.3333
$ws.Range("J80").Value = 3882.182
$ws.Range("K80").Value = 2616.3333
$ws.Range("L80").Value = 3882.182
$ws.Range("M80").Value = -1618.3333
$ws.Range("N80").Value = -5878.182
# Row 83
$ws.Range("H83").Value = 3610.9285
$ws.Range("I83").Value = 2616.3333
$ws.Range("J83").Value = 3882.182
$ws.Range("K83").Value = 13081.6665
$ws.Range("L83").Value = 19410.91
$ws.Range("M83").Value = -8089.666499999999
$ws.Range("N83").Value = -29394.91
# Row 102
$ws.Range("H102").Value = 1028.55
$ws.Range("I102").Value = 785.6875
$ws.Range("K102").Value = 785.6875
$ws.Range("M102").Value = 836.3125
# Row 126
$ws.Range("H126").Value = 4624.875
$ws.Range("I126").Value = 3499.5
$ws.Range("K126").Value = 10498.5
$ws.Range("M126").Value = -8028.5
# Row 132
$ws.Range("H132").Value = 1804.2
$ws.Range("I132").Value = 2001.75
$ws.Range("K132").Value = 6005.25
$ws.Range("M132").Value = -3475.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 33
$ws.Range("H33").Value = 28800
$ws.Range("J33").Value = 28800
$ws.Range("L33").Value = 28800
$ws.Range("N33").Value = -29380
# Row 68
$ws.Range("H68").Value = 6316.1113
$ws.Range("I68").Value = 5719
$ws.Range("J68").Value = 7062.5
$ws.Range("K68").Value = 5719
$ws.Range("L68").Value = 7062.5
$ws.Range("M68").Value = -4970
$ws.Range("N68").Value = -8560.5
# Row 71
$ws.Range("H71").Value = 6316.1113
$ws.Range("I71").Value = 5719
$ws.Range("J71").Value = 7062.5
$ws.Range("K71").Value = 28595
$ws.Range("L71").Value = 35312.5
$ws.Range("M71").Value = -24851
$ws.Range("N71").Value = -42800.5
# Row 132
$ws.Range("H132").Value = 7860
$ws.Range("I132").Value = 8200.049999999999
$ws.Range("J132").Value = 6499.8
$ws.Range("K132").Value = 24600.15
$ws.Range("L132").Value = 19499.4
$ws.Range("M132").Value = -22070.15
$ws.Range("N132").Value = -24559.4
# Row 136
$ws.Range("H136").Value = 2960.7334
$ws.Range("I136").Value = 2742.7693
$ws.Range("J136").Value = 4377.5
$ws.Range("K136").Value = 8228.3079
$ws.Range("L136").Value = 13132.5
$ws.Range("M136").Value = -5678.3079
$ws.Range("N136").Value = -18232.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
